$wb = $excel.ActiveWorkbook

# --- Update the timestamps on the existing "data" sheet (F2:F4) ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:35:42.970942"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:42.970951"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:42.970954"

# --- Add the new "metadata" sheet right after "data" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "metadata"
$newSheet.Move($null, $wb.Worksheets.Item("data"))

# Re-fetch by name: references taken before the Add/Move can go stale.
$ws = $wb.Worksheets.Item("metadata")

# Header row (B1:G1)
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$headerRange = $ws.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row (A2:G2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Short QT syndrome"
$ws.Range("C2").Value = 174
# Leading apostrophe forces this numeric-looking value to stay text (like "data_version": "0.1").
$ws.Range("D2").Value = "'0.1"
$ws.Range("E2").Value = "2020-01-15T06:46:13.208681Z"
$ws.Range("F2").Value = "2021-10-05 14:35:42.967099"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/174/?format=json"

$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$ws.Range("A1").Select() | Out-Null
